$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 65 (pushes previous row 65.. down by one),
# so the dialog-related key/value rows realign with their new content.
$ws.Rows("65:65").Insert()

# New key/value pair for the drag-force slider label.
$ws.Range("A65").Value = "dragForceSlider"
$ws.Range("B65").Value = "Drag the slider to change the Force value."

# Rows 66-74 keep their original (shifted) content: mass, accel, vel, vector,
# direction, magnitude, examples, newton_second_law_dlg_1, newton_second_law_dlg_2.
# No changes required there.

# newton_second_law_dlg_3 - reworded text.
$ws.Range("A75").Value = "newton_second_law_dlg_3"
$ws.Range("B75").Value = "As you can see, the more mass an object has, the less it will accelerate from the net force."

# newton_second_law_dlg_4 - shortened text (second sentence split off).
$ws.Range("A76").Value = "newton_second_law_dlg_4"
$ws.Range("B76").Value = "Note that: force, acceleration, and velocity are all vectors."

# New key newton_second_law_dlg_4_2 with the split-off / new explanation.
$ws.Range("A77").Value = "newton_second_law_dlg_4_2"
$ws.Range("B77").Value = "A vector is composed of a value for each axis. In our case, the x and y. These values can also tell us the direction, and the magnitude (for example: speed)."

# newton_second_law_dlg_5 now carries what used to be dlg_6's text.
$ws.Range("A78").Value = "newton_second_law_dlg_5"
$ws.Range("B78").Value = "A brave sir knight has brought us a wheel to demonstrate. The wheel is enchanted with magical trails to allow us to observe its motion."

# newton_second_law_dlg_6 now carries what used to be dlg_7's text.
$ws.Range("A79").Value = "newton_second_law_dlg_6"
$ws.Range("B79").Value = "Observe how each trail’s distance starts to increase around this area. This shows us the force being applied to the wheel."

# newton_second_law_dlg_7 - reworded text.
$ws.Range("A80").Value = "newton_second_law_dlg_7"
$ws.Range("B80").Value = "Around here, the distance between each trail is the same. This tells us that the net force on the wheel equals zero."

# newton_second_law_dlg_8 - reworded text.
$ws.Range("A81").Value = "newton_second_law_dlg_8"
$ws.Range("B81").Value = "At this point the trails are going down. We can observe that the only force acting on the wheel is the gravity."

# newton_second_law_dlg_9 now carries what used to be dlg_10's text;
# the old newton_second_law_dlg_10 key/value pair is gone.
$ws.Range("A82").Value = "newton_second_law_dlg_9"
$ws.Range("B82").Value = "The nefarious goblins have appeared out of thin air! They are surely up to no good. Get rid of them by using the wheel!"

# Update the saved selection to match the authored workbook.
$ws.Range("B65").Select()
